$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G3").Value = "mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G5").Value = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hananragab@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G6").Value = "hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
$ws.Range("G7").Value = "gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H7").Value = "154/221"
$r = $ws.Range("L10")
$r.NumberFormat = "@"
$r.Value = "43.8%"
$ws.Range("G13").Value = "wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G14").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$r = $ws.Range("S15")
$r.NumberFormat = "@"
$r.Value = "49.3%"
$r = $ws.Range("S16")
$r.NumberFormat = "@"
$r.Value = "38.2%"
$ws.Range("G24").Value = "youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G25").Value = "youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G30").Value = "hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G31").Value = "mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G32").Value = "hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G33").Value = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hananragab@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G34").Value = "hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
$ws.Range("G35").Value = "gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H35").Value = "126/246"
$ws.Range("G41").Value = "wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G42").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G52").Value = "youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G53").Value = "youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
